$d = $word.ActiveDocument

# The first paragraph (the **ID__AFFARS_...__ID** marker paragraph) gets a
# paragraph border, an updated left indent, and its id text/trailing-space
# run collapsed into a single run with the new id.
$p = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right), each with a 5pt space
# (== <w:pBdr><w:top w:space="5"/>...</w:pBdr>), no line/color/width set.
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5

# <w:ind w:left="120"/> -> <w:ind w:left="225"/>  (twips = points * 20)
$p.Format.LeftIndent = 11.25

# Replace "**ID__AFFARS_5317_topic_17__ID** " (id run + trailing space run)
# with "**ID__AFFARS_5317_703__ID**" (drops the trailing space / its run).
$rng = $p.Range
[void]$rng.Find.Execute("**ID__AFFARS_5317_topic_17__ID** ", $true, $false, $false,
                   $false, $false, $true, 1, $false,
                   "**ID__AFFARS_5317_703__ID**", 2)
